# Auto-generated edit script: apply scheduled market-data refresh to Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1569551.8
$ws.Range("I33").Value = 3448754.5
$ws.Range("J33").Value = 3549.3333
$ws.Range("K33").Value = 3448754.5
$ws.Range("L33").Value = 3549.3333
$ws.Range("M33").Value = -3448525.5
$ws.Range("N33").Value = -4007.3333
$ws.Range("H64").Value = 4499.75
$ws.Range("I64").Value = 3999.6667
$ws.Range("K64").Value = 3999.6667
$ws.Range("M64").Value = -3751.6667
$ws.Range("H67").Value = 4499.75
$ws.Range("I67").Value = 3999.6667
$ws.Range("K67").Value = 3999.6667
$ws.Range("M67").Value = -3141.6667
$ws.Range("H107").Value = 1097.375
$ws.Range("J107").Value = 1119
$ws.Range("L107").Value = 1119
$ws.Range("N107").Value = -4959
$ws.Range("H116").Value = 44188.184
$ws.Range("I116").Value = 23250
$ws.Range("K116").Value = 23250
$ws.Range("M116").Value = -19808
$ws.Range("H132").Value = 1774.037
$ws.Range("I132").Value = 1649.9615
$ws.Range("K132").Value = 4949.8845
$ws.Range("M132").Value = -2419.8845
$ws.Range("H137").Value = 3913.08
$ws.Range("J137").Value = 6188.0415
$ws.Range("L137").Value = 18564.1245
$ws.Range("N137").Value = -23664.1245

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4626549.5
$ws.Range("I32").Value = 919830
$ws.Range("K32").Value = 919830
$ws.Range("M32").Value = -919543
$ws.Range("H70").Value = 42786
$ws.Range("J70").Value = 42786
$ws.Range("L70").Value = 42786
$ws.Range("N70").Value = -43326
$ws.Range("H73").Value = 42786
$ws.Range("J73").Value = 42786
$ws.Range("L73").Value = 42786
$ws.Range("N73").Value = -44658
$ws.Range("H96").Value = 39250
$ws.Range("I96").Value = 22000
$ws.Range("K96").Value = 22000
$ws.Range("M96").Value = -19254
$ws.Range("H110").Value = 3097.4
$ws.Range("I110").Value = 2371.75
$ws.Range("J110").Value = 6000
$ws.Range("K110").Value = 2371.75
$ws.Range("L110").Value = 6000
$ws.Range("M110").Value = -326.75
$ws.Range("N110").Value = -10090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5541
$ws.Range("I20").Value = 4435.077
$ws.Range("K20").Value = 4435.077
$ws.Range("M20").Value = -4188.077
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101622
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -308112
$ws.Range("H86").Value = 1897.6666
$ws.Range("J86").Value = 2264.3333
$ws.Range("L86").Value = 2264.3333
$ws.Range("N86").Value = -4510.3333
$ws.Range("H89").Value = 1897.6666
$ws.Range("J89").Value = 2264.3333
$ws.Range("L89").Value = 11321.6665
$ws.Range("N89").Value = -22553.6665
$ws.Range("H95").Value = 26774.666
$ws.Range("J95").Value = 26774.666
$ws.Range("L95").Value = 26774.666
$ws.Range("N95").Value = -32266.666
$ws.Range("H99").Value = 32334.227
$ws.Range("I99").Value = 1541.6316
$ws.Range("J99").Value = 227354
$ws.Range("K99").Value = 1541.6316
$ws.Range("L99").Value = 227354
$ws.Range("M99").Value = -43.63159999999993
$ws.Range("N99").Value = -230350
$ws.Range("H107").Value = 2154.34
$ws.Range("I107").Value = 1998.0264
$ws.Range("J107").Value = 2649.3333
$ws.Range("K107").Value = 1998.0264
$ws.Range("L107").Value = 2649.3333
$ws.Range("M107").Value = -78.02639999999997
$ws.Range("N107").Value = -6489.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3668.5386
$ws.Range("I31").Value = 1768.125
$ws.Range("J31").Value = 6709.2
$ws.Range("K31").Value = 1768.125
$ws.Range("L31").Value = 6709.2
$ws.Range("M31").Value = -1473.125
$ws.Range("N31").Value = -7299.2
$ws.Range("H34").Value = 3668.5386
$ws.Range("I34").Value = 1768.125
$ws.Range("J34").Value = 6709.2
$ws.Range("K34").Value = 1768.125
$ws.Range("L34").Value = 6709.2
$ws.Range("M34").Value = -1566.125
$ws.Range("N34").Value = -7113.2
$ws.Range("H50").Value = 14998
$ws.Range("J50").Value = 14998
$ws.Range("L50").Value = 14998
$ws.Range("N50").Value = -16248
$ws.Range("H51").Value = 15000
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -16472
$ws.Range("H58").Value = 999.6667
$ws.Range("J58").Value = 899.5
$ws.Range("L58").Value = 899.5
$ws.Range("N58").Value = -1305.5
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 14285.571
$ws.Range("J60").Value = 15000
$ws.Range("L60").Value = 15000
$ws.Range("N60").Value = -16022
$ws.Range("H61").Value = 15000
$ws.Range("J61").Value = 15000
$ws.Range("L61").Value = 15000
$ws.Range("N61").Value = -15696
$ws.Range("H122").Value = 1539.6875
$ws.Range("I122").Value = 1333
$ws.Range("K122").Value = 3999
$ws.Range("M122").Value = -1549
$ws.Range("H136").Value = 999.6667
$ws.Range("J136").Value = 899.5
$ws.Range("L136").Value = 2698.5
$ws.Range("N136").Value = -7798.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1420.7391
$ws.Range("I68").Value = 1704.8889
$ws.Range("J68").Value = 1238.0714
$ws.Range("K68").Value = 5114.6667
$ws.Range("L68").Value = 3714.2142
$ws.Range("M68").Value = -4303.6667
$ws.Range("N68").Value = -5336.2142
$ws.Range("H71").Value = 1420.7391
$ws.Range("I71").Value = 1704.8889
$ws.Range("J71").Value = 1238.0714
$ws.Range("K71").Value = 15344.0001
$ws.Range("L71").Value = 11142.6426
$ws.Range("M71").Value = -11288.0001
$ws.Range("N71").Value = -19254.6426
$ws.Range("H107").Value = 716.2083
$ws.Range("J107").Value = 712.25
$ws.Range("L107").Value = 2136.75
$ws.Range("N107").Value = -5976.75
$ws.Range("H137").Value = 7696726
$ws.Range("J137").Value = 6587.5713
$ws.Range("L137").Value = 19762.7139
$ws.Range("N137").Value = -29962.7139
$ws.Range("H138").Value = 693.3333
$ws.Range("I138").Value = 693.3333
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 2079.9999
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 3060.0001
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12512.091
$ws.Range("I80").Value = 16670.572
$ws.Range("J80").Value = 5234.75
$ws.Range("K80").Value = 16670.572
$ws.Range("L80").Value = 5234.75
$ws.Range("M80").Value = -15672.572
$ws.Range("N80").Value = -7230.75
$ws.Range("H83").Value = 12512.091
$ws.Range("I83").Value = 16670.572
$ws.Range("J83").Value = 5234.75
$ws.Range("K83").Value = 83352.86
$ws.Range("L83").Value = 26173.75
$ws.Range("M83").Value = -78360.86
$ws.Range("N83").Value = -36157.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4428.2856
$ws.Range("I122").Value = 3250
$ws.Range("K122").Value = 9750
$ws.Range("M122").Value = -7300
$ws.Range("H136").Value = 5847.9375
$ws.Range("I136").Value = 4897.5557
$ws.Range("J136").Value = 7069.857
$ws.Range("K136").Value = 14692.6671
$ws.Range("L136").Value = 21209.571
$ws.Range("M136").Value = -12142.6671
$ws.Range("N136").Value = -26309.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 2011
$ws.Range("J20").Value = 2011
$ws.Range("L20").Value = 2011
$ws.Range("N20").Value = -2491
$ws.Range("H122").Value = 1804.317
$ws.Range("I122").Value = 1267.5518
$ws.Range("J122").Value = 3101.5
$ws.Range("K122").Value = 3802.6554
$ws.Range("L122").Value = 9304.5
$ws.Range("M122").Value = -1352.6554
$ws.Range("N122").Value = -14204.5
$ws.Range("H136").Value = 6640.5312
$ws.Range("I136").Value = 4864.16
$ws.Range("K136").Value = 14592.48
$ws.Range("M136").Value = -12042.48

Write-Output "Applied scheduled price/profit refresh to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
